$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 used to be a placeholder (date-formatted but empty) row; fill it in
# with the next week's numbers.
$ws.Range("A10").Value = 43853
$ws.Range("B10").Value = 46
$ws.Range("C10").Value = 418
$ws.Range("D10").Formula = "=C10+B10"

# Add a brand new row 11 for the following week.
$ws.Range("A11").Value = 43860
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat
$ws.Range("B11").Value = 47
$ws.Range("C11").Value = 426
$ws.Range("D11").Formula = "=C11+B11"

# Move / update the active selection to match where the author left off.
$ws.Range("B15").Select()
